$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, centered, bordered) onto the new I1/J1 headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data for rows 2-35 (col I = I0, col J = IF)
$data = @(
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 7),
    @(1, 2),
    @(1, 3),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(5, 8),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(1, 2),
    @(1, 1)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value  = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
